# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures (VALOR MORA, Cant. Trabajadores, Cant. Periodos) ---
$ws.Range("E11").Value = 92007
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3

# --- Row 16: refresh worker record (CC 1051449977 - JOSE ALFREDO CANCHILA BALLESTAS) ---
$ws.Range("C16").Value = "1051449977"
$ws.Range("D16").Value = "JOSE ALFREDO CANCHILA BALLESTAS"
$ws.Range("E16").Value = "1704"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 877803

# --- Row 17: refresh period for CC 73559072 - ELIECER ESCORCIA PADILLA ---
$ws.Range("E17").Value = "1810"
$ws.Range("G17").Value = 877803

# --- Row 20 currently carries the "last row" (bottom-border) style; load it with
#     the new final data row (CC 73559072 - ELIECER ESCORCIA PADILLA, periodo 1811)
#     before deleting the now-obsolete rows 18 & 19. Deleting those rows shifts row
#     20 up to row 18 (and the signature block from rows 25/26 up to 23/24), so the
#     table ends up with exactly the 3 data rows / correct styling in one step. ---
$ws.Range("C20").Value = "73559072"
$ws.Range("D20").Value = "ELIECER ESCORCIA PADILLA"
$ws.Range("E20").Value = "1811"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 877803

$ws.Rows("18:19").Delete()
